$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (A1:G5) has no header row -- row 1 is already a data record.
# The edit inserts a brand-new invoice record above the current row 1,
# pushing every existing row down by one (A1:G5 -> A1:G6).
$ws.Rows.Item(1).Insert()

# Stage the new record's values as TEXT in a scratch range far below the
# table (so the numeric-looking strings like "626675" aren't auto-coerced
# into numeric cells when we assign them, matching the workbook's existing
# convention of storing every value - numbers included - as a shared
# string). Using a scratch range + PasteSpecial(values) keeps the
# destination row's own (default) formatting untouched.
$scratch = $ws.Range("A1048576:G1048576")
$scratch.NumberFormat = "@"
$scratch.Cells.Item(1, 1).Value = "626675"
$scratch.Cells.Item(1, 2).Value = "Various paper supplies"
$scratch.Cells.Item(1, 3).Value = "2017-01-11"
$scratch.Cells.Item(1, 4).Value = "221966"
$scratch.Cells.Item(1, 5).Value = "44393.2"
$scratch.Cells.Item(1, 6).Value = "266359"
$scratch.Cells.Item(1, 7).Value = "RON"

$scratch.Copy()
$ws.Range("A1:G1").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$scratch.Clear()
